# Append 14 new daily OHLCV rows (rows 791-804) to the JPYUSD FX data sheet,
# extending the existing table from A1:F790 to A1:F804.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: datetime (Excel serial date), open, high, low, close, volume
$newRows = @(
    @(45142, 0.007012,  0.007063,  0.007001,  0.007053, 0),
    @(45145, 0.007053,  0.007062,  0.007014,  0.007016, 0),
    @(45146, 0.007016,  0.007021,  0.00697,   0.006971, 0),
    @(45147, 0.006974,  0.006992,  0.006954,  0.006957, 0),
    @(45148, 0.006957,  0.006974,  0.006906,  0.006908, 0),
    @(45149, 0.006907,  0.006923,  0.006895,  0.006895, 0),
    @(45152, 0.006897,  0.006908,  0.006867,  0.006867, 0),
    @(45153, 0.006869,  0.006891,  0.006856,  0.006866, 0),
    @(45154, 0.006867,  0.006881,  0.006831,  0.006832, 0),
    @(45155, 0.006833,  0.006866,  0.006826,  0.006853, 0),
    @(45156, 0.006856,  0.006898,  0.006856,  0.006874, 0),
    @(45159, 0.006874,  0.0068884, 0.006832,  0.006837, 0),
    @(45160, 0.006837,  0.006872,  0.006832,  0.006851, 0),
    @(45161, 0.0068515, 0.0068682, 0.0068515, 0.006861, 0)
)

$firstNewRow = 791
$lastExistingRow = 790
$lastNewRow = $firstNewRow + $newRows.Count - 1

# Extend the formatting (style) of the last existing data row down across the
# new rows, so that e.g. the datetime column keeps its date number format.
$ws.Range("A" + $lastExistingRow + ":F" + $lastExistingRow).Copy()
$ws.Range("A" + $firstNewRow + ":F" + $lastNewRow).PasteSpecial(-4122)

# Build a 2D array of the new values and write it in one shot.
$data = New-Object 'object[,]' $newRows.Count,6
for ($i = 0; $i -lt $newRows.Count; $i++) {
    for ($j = 0; $j -lt 6; $j++) {
        $data[$i, $j] = $newRows[$i][$j]
    }
}

$ws.Range("A" + $firstNewRow + ":F" + $lastNewRow).Value2 = $data
